$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.296321392059326
$ws.Range("B1").Value = 4.566006660461426
$ws.Range("C1").Value = 3.412265300750732
$ws.Range("D1").Value = 3.148736715316772
$ws.Range("E1").Value = 2.633186340332031
